$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:E to text format temporarily so numeric-looking strings
# (e.g. "602.92", "0.0000246") are stored as text, matching the source data,
# then restore the original style so no extra formatting is introduced.
$rng = $ws.Range("D2:E51")
$origStyle = $rng.Style
$rng.NumberFormat = "@"

$ws.Range('D2').Value = '69.204.77'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '3.789.06'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '602.92'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '164.40'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '3.788.25'
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('D11').Value = '6.32'
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '37.35'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '4.422.42'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = '3.795.73'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '69.260.71'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').Value = '7.43'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').Value = '17.38'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('E21').Value = '  +3.35%  '
$ws.Range('D22').Value = '492.51'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').Value = '0.723'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  -2.34%  '
$ws.Range('D25').Value = '84.75'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  -3.03%  '
$ws.Range('D27').Value = '12.25'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '2.98'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '8.15'
$ws.Range('E31').Value = '  +2.48%  '
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('B33').Value = 'WrappedeETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D33').Value = '3.939.87'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '31.92'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').Value = '3.735.80'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('E37').Value = '  +6.04%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').Value = '0.323'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '48.43'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').Value = '425.20'
$ws.Range('E44').Value = '  -3.04%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '1.98'
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D48').Value = '142.42'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('D49').Value = '2.814.85'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').Value = '39.88'
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('E51').Value = '  +7.51%  '

$rng.Style = $origStyle
